$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "master-ball",
    "ultra-ball",
    "great-ball",
    "poke-ball",
    "safari-ball",
    "net-ball",
    "dive-ball",
    "nest-ball",
    "repeat-ball",
    "timer-ball",
    "luxury-ball",
    "premier-ball",
    "dusk-ball",
    "heal-ball",
    "quick-ball",
    "cherish-ball",
    "smoke-ball",
    "light-ball",
    "iron-ball",
    "lure-ball",
    "level-ball",
    "moon-ball",
    "heavy-ball",
    "fast-ball",
    "friend-ball",
    "love-ball",
    "park-ball",
    "sport-ball",
    "air-balloon",
    "dream-ball",
    "snowball",
    "beast-ball",
    "left-poke-ball"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
